# Commit: Changed payments to : BACS as WorldPay requires too many changees
#
# Updates the generated/auto usernames on the active sheet (Sheet1) to the
# new naming scheme, and moves the active cell selection from D8 to E8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Manufacturer177H9_AT"
$ws.Range("A4").Value = "AuthorisedRep177H10_AT"
$ws.Range("A6").Value = "Manufacturer177H7_NU"
$ws.Range("A7").Value = "AuthorisedRep177H7_NU"

$ws.Range("E8").Select() | Out-Null
